$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and date range) ---
$hdr1 = $ws.Cells.Item(8,1)
$hdr1.Characters(21,2).Text = "17"

$hdr2 = $ws.Cells.Item(9,3)
$hdr2.Characters(27,9).Text = "4/21/2025"
$hdr2.Characters(47,9).Text = "4/27/2025"

# Row 14
$ws.Cells.Item(14,6).NumberFormat = "@"
$ws.Cells.Item(14,6).Value2 = "0"
$ws.Cells.Item(14,3).Copy()
$ws.Cells.Item(14,6).PasteSpecial(-4122)
$ws.Cells.Item(14,9).Value2 = 2
$ws.Cells.Item(14,12).Value2 = -66.666666666666
$ws.Cells.Item(14,13).Value2 = 0
$ws.Cells.Item(14,14).Value2 = -66.666666666666

# Row 15
$ws.Cells.Item(15,3).Value2 = 1
$ws.Cells.Item(15,4).Value2 = 1
$ws.Cells.Item(15,5).Value2 = 0
$ws.Cells.Item(15,6).Value2 = 4
$ws.Cells.Item(15,7).Value2 = 2
$ws.Cells.Item(15,8).Value2 = 100
$ws.Cells.Item(15,9).Value2 = 10
$ws.Cells.Item(15,10).Value2 = 8
$ws.Cells.Item(15,11).Value2 = 25
$ws.Cells.Item(15,12).Value2 = 42.857142857142
$ws.Cells.Item(15,13).Value2 = -16.666666666666
$ws.Cells.Item(15,14).Value2 = -54.545454545454

# Row 16
$ws.Cells.Item(16,3).Value2 = 1
$ws.Cells.Item(16,4).Value2 = 4
$ws.Cells.Item(16,5).Value2 = -75
$ws.Cells.Item(16,6).Value2 = 8
$ws.Cells.Item(16,7).Value2 = 15
$ws.Cells.Item(16,8).Value2 = -46.666666666666
$ws.Cells.Item(16,9).Value2 = 45
$ws.Cells.Item(16,10).Value2 = 47
$ws.Cells.Item(16,11).Value2 = -4.255319148936
$ws.Cells.Item(16,12).Value2 = -30.76923076923
$ws.Cells.Item(16,13).Value2 = -46.428571428571
$ws.Cells.Item(16,14).Value2 = -86.686390532544

# Row 17
$ws.Cells.Item(17,3).Value2 = 10
$ws.Cells.Item(17,4).Value2 = 14
$ws.Cells.Item(17,5).Value2 = -28.571428571428
$ws.Cells.Item(17,6).Value2 = 39
$ws.Cells.Item(17,7).Value2 = 41
$ws.Cells.Item(17,8).Value2 = -4.878048780487
$ws.Cells.Item(17,9).Value2 = 127
$ws.Cells.Item(17,10).Value2 = 165
$ws.Cells.Item(17,11).Value2 = -23.030303030303
$ws.Cells.Item(17,12).Value2 = 4.098360655737
$ws.Cells.Item(17,13).Value2 = 23.300970873786
$ws.Cells.Item(17,14).Value2 = -42.272727272727

# Row 18
$ws.Cells.Item(18,4).Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)
$ws.Cells.Item(18,3).Value2 = 1
$ws.Cells.Item(18,4).Value2 = 5
$ws.Cells.Item(18,5).Value2 = -80
$ws.Cells.Item(18,6).Value2 = 6
$ws.Cells.Item(18,7).Value2 = 13
$ws.Cells.Item(18,8).Value2 = -53.846153846153
$ws.Cells.Item(18,9).Value2 = 38
$ws.Cells.Item(18,10).Value2 = 37
$ws.Cells.Item(18,11).Value2 = 2.702702702702
$ws.Cells.Item(18,12).Value2 = -28.301886792452
$ws.Cells.Item(18,13).Value2 = -49.333333333333
$ws.Cells.Item(18,14).Value2 = -93.040293040293

# Row 19
$ws.Cells.Item(19,3).Value2 = 5
$ws.Cells.Item(19,4).Value2 = 5
$ws.Cells.Item(19,5).Value2 = 0
$ws.Cells.Item(19,6).Value2 = 13
$ws.Cells.Item(19,7).Value2 = 22
$ws.Cells.Item(19,8).Value2 = -40.90909090909
$ws.Cells.Item(19,9).Value2 = 81
$ws.Cells.Item(19,10).Value2 = 114
$ws.Cells.Item(19,11).Value2 = -28.947368421052
$ws.Cells.Item(19,12).Value2 = -38.636363636363
$ws.Cells.Item(19,13).Value2 = -27.678571428571
$ws.Cells.Item(19,14).Value2 = -50.306748466257

# Row 20
$ws.Cells.Item(20,4).Copy()
$ws.Cells.Item(20,3).PasteSpecial(-4122)
$ws.Cells.Item(20,3).Value2 = 2
$ws.Cells.Item(20,4).Value2 = 1
$ws.Cells.Item(20,5).Value2 = 100
$ws.Cells.Item(20,6).Value2 = 5
$ws.Cells.Item(20,7).Value2 = 5
$ws.Cells.Item(20,8).Value2 = 0
$ws.Cells.Item(20,9).Value2 = 16
$ws.Cells.Item(20,10).Value2 = 32
$ws.Cells.Item(20,11).Value2 = -50
$ws.Cells.Item(20,12).Value2 = -60
$ws.Cells.Item(20,13).Value2 = -70.37037037037
$ws.Cells.Item(20,14).Value2 = -96.428571428571

# Row 21
$ws.Cells.Item(21,3).Value2 = 20
$ws.Cells.Item(21,4).Value2 = 30
$ws.Cells.Item(21,5).Value2 = -33.333333333333
$ws.Cells.Item(21,6).Value2 = 75
$ws.Cells.Item(21,7).Value2 = 98
$ws.Cells.Item(21,8).Value2 = -23.469387755102
$ws.Cells.Item(21,9).Value2 = 319
$ws.Cells.Item(21,10).Value2 = 403
$ws.Cells.Item(21,11).Value2 = -20.843672456575
$ws.Cells.Item(21,12).Value2 = -24.941176470588
$ws.Cells.Item(21,13).Value2 = -27.828054298642
$ws.Cells.Item(21,14).Value2 = -81.698221457257

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = "0"
$ws.Cells.Item(23,3).Copy()
$ws.Cells.Item(23,4).PasteSpecial(-4122)
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value2 = "***.*"
$ws.Cells.Item(23,3).Copy()
$ws.Cells.Item(23,5).PasteSpecial(-4122)
$ws.Cells.Item(23,6).NumberFormat = "@"
$ws.Cells.Item(23,6).Value2 = "0"
$ws.Cells.Item(23,3).Copy()
$ws.Cells.Item(23,6).PasteSpecial(-4122)
$ws.Cells.Item(23,7).Value2 = 6
$ws.Cells.Item(23,8).Value2 = -100
$ws.Cells.Item(23,9).Value2 = 19
$ws.Cells.Item(23,10).Value2 = 26
$ws.Cells.Item(23,11).Value2 = -26.923076923076
$ws.Cells.Item(23,12).Value2 = -36.666666666666
$ws.Cells.Item(23,13).Value2 = 46.153846153846

# Row 24
$ws.Cells.Item(24,3).Value2 = 25
$ws.Cells.Item(24,4).Value2 = 29
$ws.Cells.Item(24,5).Value2 = -13.793103448275
$ws.Cells.Item(24,6).Value2 = 94
$ws.Cells.Item(24,7).Value2 = 102
$ws.Cells.Item(24,8).Value2 = -7.843137254901
$ws.Cells.Item(24,9).Value2 = 422
$ws.Cells.Item(24,10).Value2 = 449
$ws.Cells.Item(24,11).Value2 = -6.013363028953
$ws.Cells.Item(24,12).Value2 = 9.895833333333
$ws.Cells.Item(24,13).Value2 = 1.686746987951

# Row 25
$ws.Cells.Item(25,3).Value2 = 12
$ws.Cells.Item(25,4).Value2 = 16
$ws.Cells.Item(25,5).Value2 = -25
$ws.Cells.Item(25,6).Value2 = 44
$ws.Cells.Item(25,7).Value2 = 47
$ws.Cells.Item(25,8).Value2 = -6.382978723404
$ws.Cells.Item(25,9).Value2 = 230
$ws.Cells.Item(25,10).Value2 = 216
$ws.Cells.Item(25,11).Value2 = 6.481481481481
$ws.Cells.Item(25,12).Value2 = 60.83916083916

# Row 26
$ws.Cells.Item(26,3).Value2 = 12
$ws.Cells.Item(26,4).Value2 = 16
$ws.Cells.Item(26,5).Value2 = -25
$ws.Cells.Item(26,6).Value2 = 49
$ws.Cells.Item(26,7).Value2 = 66
$ws.Cells.Item(26,8).Value2 = -25.757575757575
$ws.Cells.Item(26,9).Value2 = 232
$ws.Cells.Item(26,10).Value2 = 214
$ws.Cells.Item(26,11).Value2 = 8.411214953271
$ws.Cells.Item(26,12).Value2 = 4.504504504504
$ws.Cells.Item(26,13).Value2 = -37.967914438502

# Row 27
$ws.Cells.Item(27,3).Value2 = 1
$ws.Cells.Item(27,4).Value2 = 1
$ws.Cells.Item(27,5).Value2 = 0
$ws.Cells.Item(27,6).Value2 = 4
$ws.Cells.Item(27,7).Value2 = 3
$ws.Cells.Item(27,8).Value2 = 33.333333333333
$ws.Cells.Item(27,9).Value2 = 17
$ws.Cells.Item(27,10).Value2 = 15
$ws.Cells.Item(27,11).Value2 = 13.333333333333
$ws.Cells.Item(27,12).Value2 = 88.888888888888

# Row 28
$ws.Cells.Item(28,3).Value2 = 2
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value2 = "0"
$ws.Cells.Item(28,13).Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4122)
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value2 = "***.*"
$ws.Cells.Item(28,13).Copy()
$ws.Cells.Item(28,5).PasteSpecial(-4122)
$ws.Cells.Item(28,6).Value2 = 11
$ws.Cells.Item(28,7).Value2 = 5
$ws.Cells.Item(28,8).Value2 = 120
$ws.Cells.Item(28,9).Value2 = 29
$ws.Cells.Item(28,10).Value2 = 28
$ws.Cells.Item(28,11).Value2 = 3.571428571428
$ws.Cells.Item(28,12).Value2 = 38.095238095238

# Row 29
$ws.Cells.Item(29,9).Copy()
$ws.Cells.Item(29,4).PasteSpecial(-4122)
$ws.Cells.Item(29,4).Value2 = 1
$ws.Cells.Item(29,12).Copy()
$ws.Cells.Item(29,5).PasteSpecial(-4122)
$ws.Cells.Item(29,5).Value2 = -100
$ws.Cells.Item(29,6).NumberFormat = "@"
$ws.Cells.Item(29,6).Value2 = "0"
$ws.Cells.Item(29,3).Copy()
$ws.Cells.Item(29,6).PasteSpecial(-4122)
$ws.Cells.Item(29,9).Copy()
$ws.Cells.Item(29,7).PasteSpecial(-4122)
$ws.Cells.Item(29,7).Value2 = 1
$ws.Cells.Item(29,12).Copy()
$ws.Cells.Item(29,8).PasteSpecial(-4122)
$ws.Cells.Item(29,8).Value2 = -100
$ws.Cells.Item(29,9).Value2 = 2
$ws.Cells.Item(29,10).Value2 = 5
$ws.Cells.Item(29,11).Value2 = -60
$ws.Cells.Item(29,12).Value2 = -50
$ws.Cells.Item(29,13).Value2 = -83.333333333333
$ws.Cells.Item(29,14).Value2 = -89.473684210526

# Row 30
$ws.Cells.Item(30,9).Copy()
$ws.Cells.Item(30,4).PasteSpecial(-4122)
$ws.Cells.Item(30,4).Value2 = 1
$ws.Cells.Item(30,12).Copy()
$ws.Cells.Item(30,5).PasteSpecial(-4122)
$ws.Cells.Item(30,5).Value2 = -100
$ws.Cells.Item(30,6).NumberFormat = "@"
$ws.Cells.Item(30,6).Value2 = "0"
$ws.Cells.Item(30,3).Copy()
$ws.Cells.Item(30,6).PasteSpecial(-4122)
$ws.Cells.Item(30,9).Copy()
$ws.Cells.Item(30,7).PasteSpecial(-4122)
$ws.Cells.Item(30,7).Value2 = 1
$ws.Cells.Item(30,12).Copy()
$ws.Cells.Item(30,8).PasteSpecial(-4122)
$ws.Cells.Item(30,8).Value2 = -100
$ws.Cells.Item(30,9).Value2 = 2
$ws.Cells.Item(30,10).Value2 = 5
$ws.Cells.Item(30,11).Value2 = -60
$ws.Cells.Item(30,12).Value2 = -50
$ws.Cells.Item(30,13).Value2 = -80
$ws.Cells.Item(30,14).Value2 = -88.235294117647

# Row 31
$ws.Cells.Item(31,10).Copy()
$ws.Cells.Item(31,6).PasteSpecial(-4122)
$ws.Cells.Item(31,6).Value2 = 1
$ws.Cells.Item(31,10).Copy()
$ws.Cells.Item(31,9).PasteSpecial(-4122)
$ws.Cells.Item(31,9).Value2 = 1
$ws.Cells.Item(31,10).Value2 = 2
$ws.Cells.Item(31,11).Value2 = -50

